# Apply the "Evaluation with inference setup" edit.
# - Evaluation sheet gets a small matrix relating datasets (rows) to models (columns).
# - Prompts sheet is restructured into labeled sections (General / Hallucination / Specific)
#   each followed by rows of prompts, with placeholder "test" answers per model column.

$wb = $excel.ActiveWorkbook

$wsEval    = $wb.Worksheets.Item("Evaluation")
$wsPrompts = $wb.Worksheets.Item("Prompts")
$wsModels  = $wb.Worksheets.Item("Models")

# A bold-styled source cell (style index 1 in the original workbook) used so that
# newly bolded cells re-use the existing style instead of creating a new one.
$boldSource = $wsModels.Range("A1")

function Set-BoldCell($ws, $addr, $text) {
    $boldSource.Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = 0
    $ws.Range($addr).Value = $text
}

# ---------------------------------------------------------------------------
# Evaluation sheet
# ---------------------------------------------------------------------------
$wsEval.Cells.Clear() | Out-Null

Set-BoldCell $wsEval "B1" "LLaMA-7b"
Set-BoldCell $wsEval "C1" "LLaMA-2-70b"
Set-BoldCell $wsEval "D1" "LLaMA-2-7b-chat"

Set-BoldCell $wsEval "A2" "I/O falsch"
$wsEval.Range("B2").Value = "guanaco-7b"

Set-BoldCell $wsEval "A3" "Alpaca"
$wsEval.Range("B3").Value = "alpaca-7b"

Set-BoldCell $wsEval "A4" "Alpaca ohne Modulenamen im Kontext"
$wsEval.Range("C4").Value = "alpaca-2-70b"

Set-BoldCell $wsEval "A5" "Alpaca refined"
$wsEval.Range("D5").Value = "alpaca-2-7b-chat"

# ---------------------------------------------------------------------------
# Prompts sheet
# ---------------------------------------------------------------------------
$wsPrompts.Cells.Clear() | Out-Null

Set-BoldCell $wsPrompts "A1" "General"

$generalPrompts = @(
    "What is a module?",
    "In what module can I edit customers?",
    "In what module do I edit the name of a customer?",
    "What is the module for entering sales invoices?"
)

$row = 2
foreach ($p in $generalPrompts) {
    $wsPrompts.Cells.Item($row, 1).Value = $p
    $wsPrompts.Cells.Item($row, 2).Value = "test"
    $wsPrompts.Cells.Item($row, 3).Value = "test"
    $wsPrompts.Cells.Item($row, 4).Value = "test"
    $wsPrompts.Cells.Item($row, 5).Value = "test"
    $row++
}

Set-BoldCell $wsPrompts "A6" "Hallucination"

$hallucinationPrompts = @(
    "This is the context of the module billofma: feeding guinea pigs and groundhogs. Which module describes feeding mammals?",
    "This is the description of the module billofma: Feeding guinea pigs and groundhogs. Which module describes Donald Trump's presidency?"
)

$row = 7
foreach ($p in $hallucinationPrompts) {
    $wsPrompts.Cells.Item($row, 1).Value = $p
    $wsPrompts.Cells.Item($row, 2).Value = "test"
    $wsPrompts.Cells.Item($row, 3).Value = "test"
    $wsPrompts.Cells.Item($row, 4).Value = "test"
    $wsPrompts.Cells.Item($row, 5).Value = "test"
    $row++
}

Set-BoldCell $wsPrompts "A9" "Specific"

Set-BoldCell $wsPrompts "A10" "balanfac"

$balanfacPrompts = @(
    "With this module, the annual and period balances of a general ledger or personal account posted in financial accounting are displayed. Which module is being described?",
    "Which module is used to display the annual and period balances of a general ledger?",
    "What is the module to list annual balances of general ledger?"
)

$row = 11
foreach ($p in $balanfacPrompts) {
    $wsPrompts.Cells.Item($row, 1).Value = $p
    $wsPrompts.Cells.Item($row, 2).Value = "test"
    $wsPrompts.Cells.Item($row, 3).Value = "test"
    $wsPrompts.Cells.Item($row, 4).Value = "test"
    $wsPrompts.Cells.Item($row, 5).Value = "test"
    $row++
}

Set-BoldCell $wsPrompts "A14" "icastedt"

$icastedtPrompts = @(
    "Which module deals with creating and deleting parts or service-role relationships?"
)

$row = 15
foreach ($p in $icastedtPrompts) {
    $wsPrompts.Cells.Item($row, 1).Value = $p
    $wsPrompts.Cells.Item($row, 2).Value = "test"
    $wsPrompts.Cells.Item($row, 3).Value = "test"
    $wsPrompts.Cells.Item($row, 4).Value = "test"
    $wsPrompts.Cells.Item($row, 5).Value = "test"
    $row++
}

Set-BoldCell $wsPrompts "A16" "billofma"

$billofmaPrompts = @(
    "Which module describes the composition of a production part?"
)

$row = 17
foreach ($p in $billofmaPrompts) {
    $wsPrompts.Cells.Item($row, 1).Value = $p
    $wsPrompts.Cells.Item($row, 2).Value = "test"
    $wsPrompts.Cells.Item($row, 3).Value = "test"
    $wsPrompts.Cells.Item($row, 4).Value = "test"
    $wsPrompts.Cells.Item($row, 5).Value = "test"
    $row++
}
